# Applies the FAQ.docx "minor changes, typos etc in prezi" edit:
#  1. "Why not RE..." - wrap "noisy" in a gramStart/gramEnd proofErr pair (no text change)
#  2. "Why not diagnostic expectations?" - append a Wingdings-arrow + Afrouzi et al sentence
#  3. "What are the key features..." - wrap trailing ")" in a gramStart/gramEnd proofErr pair
#  4. "What do agents think?..." - wrap "Susanto" in a spellStart/spellEnd proofErr pair
#  5. Insert a new bullet "Estimated changes in LR-E..." after the Susanto bullet,
#     carrying the _GoBack bookmark that used to sit on the "asymmetry" bullet
#  6. "Assumptions on Ehat?" - wrap "Ehat" in a spellStart/spellEnd proofErr pair

$d = $word.ActiveDocument

function Insert-RawXml($range, [string]$innerBodyXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

# --- 1. "Why not RE with a drift in inflation? (noisy info, Philippe Andrade)" ---
$p2 = $d.Paragraphs(4)
$xml2 = '<w:p>' + $listPPr + `
    '<w:r><w:t>Why not RE with a drift in inflation? (</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>noisy</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> info, Philippe Andrade)</w:t></w:r>' + `
    '</w:p>'
Insert-RawXml $p2.Range $xml2

# --- 2. "Why not diagnostic expectations?" gains a Wingdings arrow + Afrouzi sentence ---
$p4 = $d.Paragraphs(6)
$xml4 = '<w:p>' + $listPPr + `
    '<w:r><w:t>Why not diagnostic expectations?</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Afrouzi</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> et al 2020 WP document that diagnostic expectations don' + [char]0x2019 + 't fit the varying levels of overreaction observed in experimental data. </w:t></w:r>' + `
    '</w:p>'
Insert-RawXml $p4.Range $xml4

# --- 3. "What are the key features of the model (that distinguish it from say Phelps?)" ---
$p8 = $d.Paragraphs(8)
$xml8 = '<w:p>' + $listPPr + `
    '<w:r><w:t>What are the key features of the model (that distinguish it from say Phelps?</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>)</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
Insert-RawXml $p8.Range $xml8

# --- 4. "What do agents think? / What's in their heads? (Susanto)" ---
$p9 = $d.Paragraphs(9)
$xml9 = '<w:p>' + $listPPr + `
    '<w:r><w:t>What do agents think? / What' + [char]0x2019 + 's in their heads? (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Susanto</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>)</w:t></w:r>' + `
    '</w:p>'
Insert-RawXml $p9.Range $xml9

# --- 5. Insert new bullet after the Susanto bullet, taking over the _GoBack bookmark ---
$p9again = $d.Paragraphs(9)
$p9again.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs(10)
$xmlNew = '<w:p>' + $listPPr + `
    '<w:r><w:t>Estimated changes in LR-E: why doesn' + [char]0x2019 + 't the gain have the pattern in the motivation plot?</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'
Insert-RawXml $pNew.Range $xmlNew

# --- 6. Remove the stale bookmark from the "asymmetry" bullet (now shifted down by one) ---
$pAsym = $d.Paragraphs(13)
$xmlAsym = '<w:p>' + $listPPr + `
    '<w:r><w:t xml:space="preserve"> Why is the asymmetry not there?</w:t></w:r>' + `
    '</w:p>'
Insert-RawXml $pAsym.Range $xmlAsym

# --- 7. " Assumptions on Ehat? " ---
$pEhat = $d.Paragraphs(14)
$xmlEhat = '<w:p>' + $listPPr + `
    '<w:r><w:t xml:space="preserve"> Assumptions on </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Ehat</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">? </w:t></w:r>' + `
    '</w:p>'
Insert-RawXml $pEhat.Range $xmlEhat

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
